# Hydro implemented and tested
# Update the "Size" notation column in the HydroUnitBlock sheet to the
# standardized bracket notation used elsewhere in the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HydroUnitBlock")

$ws.Range("E8").Value  = "[NR][T]"
$ws.Range("E9").Value  = "[NA] | [T][NA]"
$ws.Range("E10").Value = "[NA] | [T][NA]"
$ws.Range("E11").Value = "[NA] | [T][NA]"
$ws.Range("E12").Value = "[NA] | [T][NA]"
$ws.Range("E13").Value = "[NA] | [T][NA]"
$ws.Range("E14").Value = "[NA] | [T][NA]"
$ws.Range("E15").Value = "[NA]"
$ws.Range("E18").Value = "[NA] | [T][NA]"
$ws.Range("E19").Value = "[NA] | [T][NA]"
$ws.Range("E22").Value = "[NA] | [T][NA]"

$ws.Activate()
$ws.Range("E15").Select()
